# Apply the BIIBNamed.xlsx update:
#  - Add a new "Down" category used by the Verdict/trend column (Y).
#  - Add a "trend" value (X) and category label (Y) for the existing row 3.
#  - Append a brand new data row (row 4) with a full set of sentiment /
#    trading metrics, matching the columns used by rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new trend score + category columns (X3, Y3) ---------------
$ws.Cells.Item(3, 24).Value = -1.3299870000000169   # X3
$ws.Cells.Item(3, 25).Value = "Down"                # Y3

# --- Row 4: brand new data row ------------------------------------------
$ws.Cells.Item(4, 1).Value = 42633.888368055559      # A4 - Date
# Reuse the date-formatted style from the row above instead of assigning a
# NumberFormat string directly (which would register a brand-new custom
# number format / cell style).
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(4, 1).PasteSpecial(-4122)
$ws.Cells.Item(4, 1).Value = 42633.888368055559
$ws.Cells.Item(4, 2).Value = 10                      # B4 - ScoreFinal
$ws.Cells.Item(4, 3).Value = "Buy"                   # C4 - Verdict
$ws.Cells.Item(4, 4).Value = 28                      # D4 - totalSentiment
$ws.Cells.Item(4, 5).Value = 16852                   # E4 - wordCount
$ws.Cells.Item(4, 6).Value = 930                     # F4 - sentenceCount
$ws.Cells.Item(4, 7).Value = 58                      # G4 - posWordPercentage
$ws.Cells.Item(4, 8).Value = 39                      # H4 - negWordPercentage
$ws.Cells.Item(4, 9).Value = 95                      # I4 - posPhrasePercentage
$ws.Cells.Item(4, 10).Value = 5                      # J4 - negPhrasePercentage
$ws.Cells.Item(4, 11).Value = 34002                  # K4 - ElapsedMs
$ws.Cells.Item(4, 12).Value = 145                    # L4 - posWordCount
$ws.Cells.Item(4, 13).Value = 99                     # M4 - negWordCount
$ws.Cells.Item(4, 14).Value = 19                     # N4 - positivePhraseCount
$ws.Cells.Item(4, 15).Value = 1                      # O4 - negativePhraseCount
$ws.Cells.Item(4, 16).Value = "Named"                # P4 - Method
$ws.Cells.Item(4, 17).Value = 0                      # Q4 - RSI
$ws.Cells.Item(4, 18).Value = 1.76                   # R4 - PEG
$ws.Cells.Item(3, 19).Copy()
$ws.Cells.Item(4, 19).PasteSpecial(-4122)
$ws.Cells.Item(4, 19).Value = 0.1055                 # S4 - 200Moving%
$ws.Cells.Item(4, 20).Value = -6.67                  # T4 - 50Moving%
$ws.Cells.Item(4, 21).Value = 5.83                   # U4 - PriceBook
$ws.Cells.Item(4, 22).Value = "N/A"                  # V4 - Dividend
$ws.Cells.Item(4, 23).Value = 0                      # W4 - Bollinger
